$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# 英語 -> 英语 (both occurrences are identical replacements)
Replace-Text "英語" "英语"

# / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語 -> / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语
# (leading space intentionally excluded from the match to avoid Word's Find/Replace
#  bleeding the preceding hyperlink run's formatting into this run)
Replace-Text "葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" "葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

# 簡介 -> 简介
Replace-Text "簡介" "简介"

# Intro paragraph
Replace-Text "發送給目標國家中那些文件未通過我們驗證流程的合作夥伴的電子郵件。 將通過 customer.io 發送" "一封发送给目标国家中未通过我们验证流程的合作伙伴的电子邮件。 将通过 customer.io 发送"

# 目標受眾 -> 目标受众
Replace-Text "目標受眾" "目标受众"

# Target audience description
Replace-Text "提交錯誤/不完整文檔的被邀請合作夥伴" "提交了错误/不完整文件的邀请合作伙伴"

# 主題行 -> 主题行
Replace-Text "主題行" "主题行"

# [事件名稱] -> [事件名称]
Replace-Text "[事件名稱]" "[事件名称]"

# — 文件驗證失敗 -> — 文档验证失败
Replace-Text " — 文件驗證失敗 " " — 文档验证失败 "

# 啊哦！ 文檔無法驗證 -> 啊哦！ 文件无法验证
Replace-Text "啊哦！ 文檔無法驗證" "啊哦！ 文件无法验证"

# [合作夥伴姓名] -> [合作伙伴姓名]
Replace-Text "[合作夥伴姓名]" "[合作伙伴姓名]"

# Regret sentence -> English
Replace-Text "很遺憾地通知您，您的文檔未通過驗證流程，因為我們發現以下問題： " "We regret to inform you that your documents have failed our verification process as we found the following issues with them: "

# 您的疫苗接種證明副本 -> 疫苗接种证书副本
Replace-Text "您的疫苗接種證明副本" "疫苗接种证书副本"

# : 文檔不清楚 -> : 文件不清楚
Replace-Text ": 文檔不清楚" ": 文件不清楚"

# [文檔 2] -> [文件 2]
Replace-Text "[文檔 2]" "[文件 2]"

# : [問題] -> : [problem]
Replace-Text ": [問題]" ": [problem]"

# 請在  -> 请在  (trailing space preserved)
Replace-Text "請在 " "请在 "

# trailing: 之前重新提交上述文檔，以便我們進行必要的安排。
Replace-Text " 之前重新提交上述文檔，以便我們進行必要的安排。" " 之前重新提交上述文件，以便我们进行必要的安排。"

# 如有任何疑問，請通過  -> 如有任何疑问，请通过 
Replace-Text "如有任何疑問，請通過 " "如有任何疑问，请通过 "

# [電子郵件地址] -> [电子邮件地址]
Replace-Text "[電子郵件地址]" "[电子邮件地址]"

# [WHATSAPP 號碼] -> [WHATSAPP 号码]
Replace-Text "[WHATSAPP 號碼]" "[WHATSAPP 号码]"

# (WhatsApp) 聯繫您的區域經理,  -> (WhatsApp) 联系您的区域经理  (comma removed)
Replace-Text " (WhatsApp) 聯繫您的區域經理, " " (WhatsApp) 联系您的区域经理 "

# [姓名] -> [NAME]
Replace-Text "[姓名]" "[NAME]"

# 。  ->  。  (leading space added) - scope the search to the text after [NAME]
# so we don't touch the unrelated "。" earlier in the document (inside the
# already-translated intro paragraph). Note: ReplaceAll (2) on a sub-range
# still replaces matches in the *whole* document in this runtime, so we use
# wdReplaceOne (1) which correctly honours the range/start position and only
# replaces the single nearest match.
$nameRange = $d.Content
$nameRange.Find.Execute("[NAME]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($nameRange.Find.Found) {
    $afterName = $d.Range($nameRange.End, $d.Content.End)
    $afterName.Find.Execute("。 ", $true, $false, $false, $false, $false, $true, 1, $false, " 。 ", 1) | Out-Null
}
